# Daily attendance processing - 2025-12-19 11:27:56
# Reverses the order of comma-separated "Recorded By" entries in column G
# for specific known value combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

$map = @{
    "System, system, backup@backdoor.com" = "backup@backdoor.com, system, System"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
